$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: title changed
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 28: title + link changed
$ws.Range("D28").Value = "좋은 그리퍼는 파지 견고성이 높다"
$ws.Range("E28").Value = "https://ropiens.tistory.com/188"

# Row 36: title + link changed
$ws.Range("D36").Value = "Change point detection in time series"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/362"

# Row 51: title + link changed
$ws.Range("D51").Value = "[python+pandas] groupby 메소드로 그룹의 평균값, 최대값 산출하기"
$ws.Range("E51").Value = "https://bskyvision.com/978"
